# eims-toi-transect/toi-transect.xlsx
# "get toi_source and API data for ncp-gop from eims-toi output"
#
# 1. Insert a new row at row 2 for the "cruise" attribute (pushes the
#    existing attribute rows down by one).
# 2. Fill in the new row 2 with the cruise attribute metadata.
# 3. Re-order / rename: what is now row 10 ("depth_matlab", shifted down
#    from old row 9) and row 11 ("depth_API", shifted down from old row
#    10) need to become, in order, "depth" (renamed from depth_API, with
#    an updated definition) followed by "depth_matlab".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1 & selection: insert new row 2, leaving the freshly-inserted row
# selected the way Excel does after an Insert-Row command.
$ws.Rows(2).Insert()
$ws.Rows(2).Select()

# --- 2: populate the new "cruise" attribute row.
$ws.Range("A2").Value = "cruise"
$ws.Range("B2").Value = "Identifier for research cruise generally including abbreviation for research vessel and voyage number"
$ws.Range("C2").Value = "character"
# The inserted row copied formatting (wrap-text) down from row 1's B
# column; the source row had no such style, so strip it back off.
$ws.Range("B2").ClearFormats()

# --- 3: rewrite rows 10-11 into their final order/content.
$ws.Range("A10").Value = "depth"
$ws.Range("B10").Value = "Data product depth of sample below sea surface, for underway samples depth of ship's intake, for Niskins from CTD summary data in NES-LTER  API"
$ws.Range("C10").Value = "numeric"
$ws.Range("D10").Value = "meter"

$ws.Range("A11").Value = "depth_matlab"
$ws.Range("B11").Value = "PI-provided depth of sample below sea surface. URI http://vocab.nerc.ac.uk/collection/P09/current/DEPH/"
$ws.Range("C11").Value = "numeric"
$ws.Range("D11").Value = "meter"
